$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: 204. Count Primes -------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C33"), "https://leetcode.com/problems/count-primes/")
$ws.Range("A29:H29").Copy()
$ws.Range("A33:H33").PasteSpecial(-4122)

$ws.Range("A33").Value = "204. Count Primes"
$ws.Range("B33").Value = "Easy"
$ws.Range("C33").Value = "https://leetcode.com/problems/count-primes/"
$ws.Range("D33").Value = 44536
$ws.Range("E33").Value = "质数"
$ws.Range("F33").Value = "按常规的一个个判断会超时，从奇数、开平方、倍数角度加快速度"
$ws.Range("G33").Value = "未复习"
$ws.Range("H33").Value = "⭕"

$ws.Range("G33").Font.Name = "宋体"
$ws.Range("G33").Font.Family = 3
$ws.Range("G33").NumberFormat = "mm-dd-yy"

$ws.Rows.Item(33).RowHeight = 28

# --- Row 34: 241. Different Ways to Add Parentheses ---------------------------
$ws.Hyperlinks.Add($ws.Range("C34"), "https://leetcode.com/problems/different-ways-to-add-parentheses/")
$ws.Range("A29:H29").Copy()
$ws.Range("A34:H34").PasteSpecial(-4122)

$ws.Range("A34").Value = "241. Different Ways to Add Parentheses"
$ws.Range("B34").Value = "Medium"
$ws.Range("C34").Value = "https://leetcode.com/problems/different-ways-to-add-parentheses/"
$ws.Range("D34").Value = 44537
$ws.Range("E34").Value = "分治"
$ws.Range("F34").Value = "按操作符的位置进行分割dfs"
$ws.Range("G34").Value = "未复习"
$ws.Range("H34").Value = "⭕"

$ws.Range("G34").Font.Name = "宋体"
$ws.Range("G34").Font.Family = 3
$ws.Range("G34").NumberFormat = "mm-dd-yy"

$ws.Rows.Item(34).RowHeight = 42

# --- View state -----------------------------------------------------------------
$ws.Range("F30").Select()
